$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C:F hold numeric-looking values that must stay text cells (as in
# the source data), so force Text format before writing the values.
$ws.Range("C2:F14").NumberFormat = "@"

# New ordering/content for rows 2-14 (columns C:F = runs, balls, fours, sixes)
# Row 14 is a brand-new row; the rest are the same underlying records,
# reordered, per the commit "updated activity till excel form".
$data = @(
    @("Faf du Plessis ", "Chennai Super Kings", "1", "7", "0", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "48", "34", "4", "2"),
    @("Faf du Plessis ", "Chennai Super Kings", "25", "13", "2", "2"),
    @("Faf du Plessis ", "Chennai Super Kings", "87", "53", "11", "1"),
    @("Faf du Plessis ", "Chennai Super Kings", "17", "10", "3", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "72", "37", "1", "7"),
    @("Faf du Plessis ", "Chennai Super Kings", "10", "9", "1", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "0", "1", "0", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "58", "47", "6", "2"),
    @("Faf du Plessis ", "Chennai Super Kings", "22", "19", "4", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "58", "44", "6", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "8", "10", "0", "0"),
    @("Faf du Plessis ", "Chennai Super Kings", "43", "35", "4", "0")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
